$wb = $excel.ActiveWorkbook
$excel.Calculation = -4135  # xlCalculationManual: preserve stale cached formula results (e.g. #DIV/0!) matching original edit

$ws = $wb.Worksheets.Item("FS")
$ws.Range("C24").Value = 1.0
$ws.Range("D24").Value = 0.0
$ws.Range("E24").Value = 0.0
$ws.Range("G24").Value = 36.84210526315789
$ws.Range("J24").Value = 0.7777777777777778
$ws.Range("K24").Value = 0.2222222222222222
$ws.Range("L24").Value = 0.0
$ws.Range("N24").Value = 31.57894736842105
$ws.Range("C25").Value = 0.125
$ws.Range("D25").Value = 0.5
$ws.Range("E25").Value = 0.375
$ws.Range("J25").Value = 0.09722222222222222
$ws.Range("K25").Value = 0.6111111111111112
$ws.Range("L25").Value = 0.2916666666666667
$ws.Range("C26").Value = 0.0
$ws.Range("D26").Value = 0.5
$ws.Range("E26").Value = 0.5
$ws.Range("J26").Value = 0.018518518518518517
$ws.Range("K26").Value = 0.2777777777777778
$ws.Range("L26").Value = 0.7037037037037037

$ws = $wb.Worksheets.Item("IF")
$ws.Range("C24").Value = 0.6
$ws.Range("D24").Value = 0.4
$ws.Range("E24").Value = 0.0
$ws.Range("G24").Value = 21.052631578947366
$ws.Range("J24").Value = 0.7333333333333333
$ws.Range("K24").Value = 0.17777777777777778
$ws.Range("L24").Value = 0.08888888888888889
$ws.Range("N24").Value = 29.82456140350877
$ws.Range("C25").Value = 0.0
$ws.Range("D25").Value = 1.0
$ws.Range("E25").Value = 0.0
$ws.Range("J25").Value = 0.08333333333333333
$ws.Range("K25").Value = 0.7361111111111112
$ws.Range("L25").Value = 0.18055555555555555
$ws.Range("C26").Value = 0.16666666666666666
$ws.Range("D26").Value = 0.16666666666666666
$ws.Range("E26").Value = 0.6666666666666666
$ws.Range("J26").Value = 0.07407407407407407
$ws.Range("K26").Value = 0.2962962962962963
$ws.Range("L26").Value = 0.6296296296296297

$ws = $wb.Worksheets.Item("IA")
$ws.Range("C24").Value = 1.0
$ws.Range("D24").Value = 0.0
$ws.Range("E24").Value = 0.0
$ws.Range("G24").Value = 21.052631578947366
$ws.Range("J24").Value = 0.8222222222222222
$ws.Range("K24").Value = 0.08888888888888889
$ws.Range("L24").Value = 0.08888888888888889
$ws.Range("N24").Value = 32.16374269005848
$ws.Range("C25").Value = 0.0
$ws.Range("D25").Value = 0.875
$ws.Range("E25").Value = 0.125
$ws.Range("J25").Value = 0.05555555555555555
$ws.Range("K25").Value = 0.6805555555555556
$ws.Range("L25").Value = 0.2638888888888889
$ws.Range("C26").Value = 0.0
$ws.Range("D26").Value = 0.5
$ws.Range("E26").Value = 0.5
$ws.Range("J26").Value = 0.05555555555555555
$ws.Range("K26").Value = 0.3888888888888889
$ws.Range("L26").Value = 0.5555555555555556

$ws = $wb.Worksheets.Item("FS-IF")
$ws.Range("C24").Value = 0.8
$ws.Range("D24").Value = 0.0
$ws.Range("E24").Value = 0.2
$ws.Range("G24").Value = 21.052631578947366
$ws.Range("J24").Value = 0.8666666666666667
$ws.Range("K24").Value = 0.13333333333333333
$ws.Range("L24").Value = 0.0
$ws.Range("N24").Value = 22.22222222222222
$ws.Range("C25").Value = 0.0
$ws.Range("D25").Value = 1.0
$ws.Range("E25").Value = 0.0
$ws.Range("J25").Value = 0.1111111111111111
$ws.Range("K25").Value = 0.7638888888888888
$ws.Range("L25").Value = 0.125
$ws.Range("C26").Value = 0.16666666666666666
$ws.Range("D26").Value = 0.3333333333333333
$ws.Range("E26").Value = 0.5
$ws.Range("J26").Value = 0.018518518518518517
$ws.Range("K26").Value = 0.25925925925925924
$ws.Range("L26").Value = 0.7222222222222222

$ws = $wb.Worksheets.Item("FS-IA")
$ws.Range("C24").Value = 1.0
$ws.Range("D24").Value = 0.0
$ws.Range("E24").Value = 0.0
$ws.Range("G24").Value = 26.31578947368421
$ws.Range("J24").Value = 0.8222222222222222
$ws.Range("K24").Value = 0.15555555555555556
$ws.Range("L24").Value = 0.022222222222222223
$ws.Range("N24").Value = 29.239766081871345
$ws.Range("C25").Value = 0.125
$ws.Range("D25").Value = 0.625
$ws.Range("E25").Value = 0.25
$ws.Range("J25").Value = 0.1111111111111111
$ws.Range("K25").Value = 0.6666666666666666
$ws.Range("L25").Value = 0.2222222222222222
$ws.Range("C26").Value = 0.0
$ws.Range("D26").Value = 0.3333333333333333
$ws.Range("E26").Value = 0.6666666666666666
$ws.Range("J26").Value = 0.018518518518518517
$ws.Range("K26").Value = 0.3148148148148148
$ws.Range("L26").Value = 0.6666666666666666

$ws = $wb.Worksheets.Item("IF-IA")
$ws.Range("C24").Value = 0.8
$ws.Range("D24").Value = 0.2
$ws.Range("E24").Value = 0.0
$ws.Range("G24").Value = 21.052631578947366
$ws.Range("J24").Value = 0.8444444444444444
$ws.Range("K24").Value = 0.13333333333333333
$ws.Range("L24").Value = 0.022222222222222223
$ws.Range("N24").Value = 24.561403508771928
$ws.Range("C25").Value = 0.0
$ws.Range("D25").Value = 0.75
$ws.Range("E25").Value = 0.25
$ws.Range("J25").Value = 0.041666666666666664
$ws.Range("K25").Value = 0.75
$ws.Range("L25").Value = 0.20833333333333334
$ws.Range("C26").Value = 0.0
$ws.Range("D26").Value = 0.16666666666666666
$ws.Range("E26").Value = 0.8333333333333334
$ws.Range("J26").Value = 0.018518518518518517
$ws.Range("K26").Value = 0.2962962962962963
$ws.Range("L26").Value = 0.6851851851851852

$ws = $wb.Worksheets.Item("FS-IF-IA")
$ws.Range("C18").Value = 0.6
$ws.Range("D18").Value = 0.4
$ws.Range("E18").Value = 0.0
$ws.Range("G18").Value = 26.31578947368421
$ws.Range("J18").Value = 0.6888888888888889
$ws.Range("K18").Value = 0.28888888888888886
$ws.Range("L18").Value = 0.022222222222222223
$ws.Range("N18").Value = 22.807017543859647
$ws.Range("C19").Value = 0.125
$ws.Range("D19").Value = 0.75
$ws.Range("E19").Value = 0.125
$ws.Range("J19").Value = 0.125
$ws.Range("K19").Value = 0.7361111111111112
$ws.Range("L19").Value = 0.1388888888888889
$ws.Range("C20").Value = 0.0
$ws.Range("D20").Value = 0.16666666666666666
$ws.Range("E20").Value = 0.8333333333333334
$ws.Range("J20").Value = 0.0
$ws.Range("K20").Value = 0.1111111111111111
$ws.Range("L20").Value = 0.8888888888888888
$ws.Range("C24").Value = 0.8
$ws.Range("D24").Value = 0.0
$ws.Range("E24").Value = 0.2
$ws.Range("G24").Value = 21.052631578947366
$ws.Range("J24").Value = 0.8666666666666667
$ws.Range("K24").Value = 0.13333333333333333
$ws.Range("L24").Value = 0.0
$ws.Range("N24").Value = 24.561403508771928
$ws.Range("C25").Value = 0.0
$ws.Range("D25").Value = 0.75
$ws.Range("E25").Value = 0.25
$ws.Range("J25").Value = 0.06944444444444445
$ws.Range("K25").Value = 0.6944444444444444
$ws.Range("L25").Value = 0.2361111111111111
$ws.Range("C26").Value = 0.0
$ws.Range("D26").Value = 0.16666666666666666
$ws.Range("E26").Value = 0.8333333333333334
$ws.Range("J26").Value = 0.0
$ws.Range("K26").Value = 0.25925925925925924
$ws.Range("L26").Value = 0.7407407407407407
